$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column F ("Is Significant"), shifting it to H.
$ws.Range("F:G").Insert()

# New header cells
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"

# Copy formatting from an existing header cell (E1) to the new headers
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill Observed / Expected data for rows 2 and 3
$ws.Range("F2").Value = "[768 188] ; [707  45]"
$ws.Range("G2").Value = "[825.58548009 130.41451991] ; [649.41451991 102.58548009]"

$ws.Range("F3").Value = "[734 245] ; [632 122]"
$ws.Range("G3").Value = "[771.67570687 207.32429313] ; [594.32429313 159.67570687]"
